# Remove the "Sunny/Rainy Test case" column (F) contents from the RTM sheet.
# F1 keeps its header style but loses its text; F4..F18 (the per-row
# Sunny/Rainy markers) are cleared entirely so the cells disappear from the
# sheet XML and the now-unused shared strings ("Sunny/Rainy Test case",
# "Sunny", "Rainy") drop out of sharedStrings.xml.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1:F18").ClearContents()

# Move/restore the active selection to I8, matching the saved selection
# state in the workbook.
$ws.Range("I8").Select()
